$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph and the "(c) 2020 ..." paragraph
# by their text content (robust against any pre-existing paragraph
# numbering assumptions).
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($startPara -eq $null -and $t -like "*Ver no Jupiter*") {
        $startPara = $i
    }
    if ($t -like "*Creative Commons Attribution*") {
        $endPara = $i
    }
}

# The blank paragraph immediately preceding "Ver no Jupiter ..." is also
# removed, while the blank paragraph that follows the copyright notice is
# kept (it stays right before the page-break paragraph).
$deleteFromPara = $startPara - 1

$deleteStart = $d.Paragraphs.Item($deleteFromPara).Range.Start
$deleteEnd = $d.Paragraphs.Item($endPara).Range.End

$r = $d.Range($deleteStart, $deleteEnd)
$r.Delete()
